$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the contents of row 4 and row 5 for columns A, B, D, E, F, G, H, Q, R
# (columns whose values differ between the two records).
$cols = @("A", "B", "D", "E", "F", "G", "H", "Q", "R")

foreach ($col in $cols) {
    $cell4 = $ws.Range($col + "4")
    $cell5 = $ws.Range($col + "5")

    $val4 = $cell4.Value2
    $val5 = $cell5.Value2

    $cell4.Value2 = $val5
    $cell5.Value2 = $val4
}
